$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.782637532668057
$ws.Range("D2").Value = 4.130189461739612
$ws.Range("E2").Value = 16.52654179159011
$ws.Range("F2").Value = 22.33544634103874
$ws.Range("G2").Value = 3.608341058757886
$ws.Range("K2").Value = 11.08149156305305
$ws.Range("N2").Value = 17.59696812025692
$ws.Range("O2").Value = 19.7886948003225

$ws.Range("B3").Value = 7.709626176934995
$ws.Range("D3").Value = 4.095399305861559
$ws.Range("E3").Value = 15.58699721510027
$ws.Range("F3").Value = 22.26299637456494
$ws.Range("G3").Value = 3.610525931659404
$ws.Range("K3").Value = 10.55227450858571
$ws.Range("N3").Value = 17.65680227979854
$ws.Range("O3").Value = 19.80072070309436

$ws.Range("B4").Value = 7.66628926026857
$ws.Range("D4").Value = 4.073585185786377
$ws.Range("E4").Value = 14.98516023932362
$ws.Range("F4").Value = 22.22593950654022
$ws.Range("G4").Value = 3.611937223591782
$ws.Range("K4").Value = 10.21127959792869
$ws.Range("N4").Value = 17.69538809553978
$ws.Range("O4").Value = 19.81393917279548

$ws.Range("B5").Value = 7.649024257641219
$ws.Range("D5").Value = 4.064585785037426
$ws.Range("E5").Value = 14.7339012490767
$ws.Range("F5").Value = 22.21271560735728
$ws.Range("G5").Value = 3.612529941715195
$ws.Range("K5").Value = 10.06838083233162
$ws.Range("N5").Value = 17.71157799854197
$ws.Range("O5").Value = 19.82078864353196

$ws.Range("B6").Value = 7.646181849195955
$ws.Range("D6").Value = 4.063084914304797
$ws.Range("E6").Value = 14.69182607629404
$ws.Range("F6").Value = 22.21063337945532
$ws.Range("G6").Value = 3.612629427202085
$ws.Range("K6").Value = 10.04441775160321
$ws.Range("N6").Value = 17.71429449613002
$ws.Range("O6").Value = 19.82201420140297

$ws.Range("B7").Value = 7.666054792763623
$ws.Range("D7").Value = 4.073464256823239
$ws.Range("E7").Value = 14.98179559437641
$ws.Range("F7").Value = 22.22575355396352
$ws.Range("G7").Value = 3.61194514583239
$ws.Range("K7").Value = 10.20936823069637
$ws.Range("N7").Value = 17.69560455019635
$ws.Range("O7").Value = 19.81402563060763

$ws.Range("B8").Value = 7.757165623521593
$ws.Range("D8").Value = 4.118290118392514
$ws.Range("E8").Value = 16.20791151230536
$ws.Range("F8").Value = 22.3089303710828
$ws.Range("G8").Value = 3.60907995853782
$ws.Range("K8").Value = 10.90240200941654
$ws.Range("N8").Value = 17.61721630954608
$ws.Range("O8").Value = 19.79162805832296

$ws.Range("B9").Value = 7.946768863479892
$ws.Range("D9").Value = 4.20243335136307
$ws.Range("E9").Value = 18.50108184575906
$ws.Range("F9").Value = 22.53044123864591
$ws.Range("G9").Value = 3.604012165563931
$ws.Range("K9").Value = 12.13063882264533
$ws.Range("N9").Value = 17.47809472315244
$ws.Range("O9").Value = 19.79415351386078

$ws.Range("B10").Value = 8.091487025110531
$ws.Range("D10").Value = 4.261743194632605
$ws.Range("E10").Value = 20.14690357292098
$ws.Range("F10").Value = 22.72789653665683
$ws.Range("G10").Value = 3.600620774166785
$ws.Range("K10").Value = 12.94965427811702
$ws.Range("N10").Value = 17.38469569216099
$ws.Range("O10").Value = 19.8244746395646

$ws.Range("B11").Value = 8.158226950868196
$ws.Range("D11").Value = 4.288134941749049
$ws.Range("E11").Value = 20.85339670585262
$ws.Range("F11").Value = 22.82501835765207
$ws.Range("G11").Value = 3.599149186199192
$ws.Range("K11").Value = 13.30360008877161
$ws.Range("N11").Value = 17.34410179252936
$ws.Range("O11").Value = 19.84446108782143

$ws.Range("B12").Value = 8.183607075612013
$ws.Range("D12").Value = 4.29804047605843
$ws.Range("E12").Value = 21.11489575662811
$ws.Range("F12").Value = 22.86282125739863
$ws.Range("G12").Value = 3.598602105548148
$ws.Range("K12").Value = 13.43491743341628
$ws.Range("N12").Value = 17.32900093274502
$ws.Range("O12").Value = 19.85291917783449

$ws.Range("B13").Value = 8.17813661461752
$ws.Range("D13").Value = 4.295911135361726
$ws.Range("E13").Value = 21.05884511254096
$ws.Range("F13").Value = 22.85463455437178
$ws.Range("G13").Value = 3.59871947739699
$ws.Range("K13").Value = 13.40675709852061
$ws.Range("N13").Value = 17.33224113097534
$ws.Range("O13").Value = 19.85105803505283

$ws.Range("B14").Value = 8.160312974112074
$ws.Range("D14").Value = 4.28895167011619
$ws.Range("E14").Value = 20.87503113981786
$ws.Range("F14").Value = 22.82810803125664
$ws.Range("G14").Value = 3.599103973890935
$ws.Range("K14").Value = 13.31445822070801
$ws.Range("N14").Value = 17.34285400763075
$ws.Range("O14").Value = 19.84513911899469

$ws.Range("B15").Value = 8.149408747059073
$ws.Range("D15").Value = 4.284677166809076
$ws.Range("E15").Value = 20.76165489058002
$ws.Range("F15").Value = 22.8119924862328
$ws.Range("G15").Value = 3.599340812862025
$ws.Range("K15").Value = 13.25756808939565
$ws.Range("N15").Value = 17.34938998124679
$ws.Range("O15").Value = 19.84162942602791

$ws.Range("B16").Value = 8.087141674248278
$ws.Range("D16").Value = 4.260006260162593
$ws.Range("E16").Value = 20.09988630735643
$ws.Range("F16").Value = 22.72169421394389
$ws.Range("G16").Value = 3.600718373219571
$ws.Range("K16").Value = 12.92614502835827
$ws.Range("N16").Value = 17.38738660764368
$ws.Range("O16").Value = 19.82329304694317

$ws.Range("B17").Value = 8.049158562774579
$ws.Range("D17").Value = 4.244717965411429
$ws.Range("E17").Value = 19.68312845239523
$ws.Range("F17").Value = 22.6681517072052
$ws.Range("G17").Value = 3.601581650461239
$ws.Range("K17").Value = 12.71802883638894
$ws.Range("N17").Value = 17.41118052395456
$ws.Range("O17").Value = 19.81363018426291

$ws.Range("B18").Value = 8.027398435220132
$ws.Range("D18").Value = 4.235869362987929
$ws.Range("E18").Value = 19.43944634971242
$ws.Range("F18").Value = 22.63804404945243
$ws.Range("G18").Value = 3.602084887210472
$ws.Range("K18").Value = 12.5965741115844
$ws.Range("N18").Value = 17.42504447543677
$ws.Range("O18").Value = 19.80865541805899

$ws.Range("B19").Value = 8.020046439384343
$ws.Range("D19").Value = 4.232864018651465
$ws.Range("E19").Value = 19.35625620772385
$ws.Range("F19").Value = 22.62796908833924
$ws.Range("G19").Value = 3.602256427426843
$ws.Range("K19").Value = 12.55515196216346
$ws.Range("N19").Value = 17.42976923337745
$ws.Range("O19").Value = 19.80707119908213

$ws.Range("B20").Value = 8.053193119830754
$ws.Range("D20").Value = 4.246351169232693
$ws.Range("E20").Value = 19.72790408791418
$ws.Range("F20").Value = 22.67378030399727
$ws.Range("G20").Value = 3.601489059852396
$ws.Range("K20").Value = 12.74036474542276
$ws.Range("N20").Value = 17.40862917388367
$ws.Range("O20").Value = 19.8145984715365

$ws.Range("B21").Value = 8.165545485795477
$ws.Range("D21").Value = 4.29099826407838
$ws.Range("E21").Value = 20.92918526855625
$ws.Range("F21").Value = 22.83587190295759
$ws.Range("G21").Value = 3.598990762179606
$ws.Range("K21").Value = 13.34164254909762
$ws.Range("N21").Value = 17.33972939816544
$ws.Range("O21").Value = 19.84685351570232

$ws.Range("B22").Value = 8.23958819355172
$ws.Range("D22").Value = 4.319660482818801
$ws.Range("E22").Value = 21.67913335804562
$ws.Range("F22").Value = 22.94776957632405
$ws.Range("G22").Value = 3.597417277956439
$ws.Range("K22").Value = 13.71878009587396
$ws.Range("N22").Value = 17.29627946893336
$ws.Range("O22").Value = 19.87311847132974

$ws.Range("B23").Value = 8.200021680570833
$ws.Range("D23").Value = 4.304411478047277
$ws.Range("E23").Value = 21.28207789980837
$ws.Range("F23").Value = 22.88751087170661
$ws.Range("G23").Value = 3.59825166904152
$ws.Range("K23").Value = 13.51895322471954
$ws.Range("N23").Value = 17.31932531564457
$ws.Range("O23").Value = 19.85862657362876

$ws.Range("B24").Value = 8.051368855382453
$ws.Range("D24").Value = 4.245612981729533
$ws.Range("E24").Value = 19.70767374232158
$ws.Range("F24").Value = 22.67123351197534
$ws.Range("G24").Value = 3.601530898483975
$ws.Range("K24").Value = 12.73027230344669
$ws.Range("N24").Value = 17.40978206438298
$ws.Range("O24").Value = 19.81415890033742

$ws.Range("B25").Value = 7.894433015970831
$ws.Range("D25").Value = 4.180095688618747
$ws.Range("E25").Value = 17.85734778435182
$ws.Range("F25").Value = 22.46434628144703
$ws.Range("G25").Value = 3.605324570467196
$ws.Range("K25").Value = 11.81279227956305
$ws.Range("N25").Value = 17.51417693264564
$ws.Range("O25").Value = 19.81402563060763
